$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the caveat text for the "achievements" rate row (row 15) and the
# "participation" rate row (row 16) so the per-100k rate wording is
# clarified. This also orphans (and therefore removes) the old shared
# string that previously only row 15 referenced, which is the same
# compaction the original commit's diff shows in sharedStrings.xml.
$ws.Range("C15").Value = "The rates are the number of achievments in AY21/22 per 100,000 of the population. Further education and skills include all age apprenticeships and publicly-funded adult (19+) learning, including community learning, delivered by an FE institution, a training provider or within a local community. `n"
$ws.Range("C16").Value = "The rates are the number of participants in AY21/22 per 100,000 of the population. Further education and skills include all age apprenticeships and publicly-funded adult (19+) learning, including community learning, delivered by an FE institution, a training provider or within a local community. `n"
